$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.532.22'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '1.882.22'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.79%  '
$ws.Range("E6").Value = '  +2.45%  '
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.95'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.35%  '
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0705'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.98%  '
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("E12").Value = '  +1.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '12.32'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.75%  '
$ws.Range("D15").Value = '1.863.44'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("E16").Value = '  +2.77%  '
$ws.Range("D17").Value = '35.543.92'
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '71.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("E19").Value = '  +2.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '243.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("E21").Value = '  +1.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.96%  '
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("E24").Value = '  +1.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.76'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("E26").Value = '  +25.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.90'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.53%  '
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0565'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.28%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.00%  '
$ws.Range("B32").Value = 'BinanceUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.939'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +24.98%  '
$ws.Range("E34").Value = '  +3.61%  '
$ws.Range("E35").Value = '  +10.90%  '
$ws.Range("E36").Value = '  +5.37%  '
$ws.Range("E37").Value = '  +10.94%  '
$ws.Range("E38").Value = '  +3.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0204'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '90.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("D41").Value = '1.354.40'
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.87%  '
$ws.Range("B43").Value = 'Gas'
$ws.Range("C43").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +48.86%  '
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0594'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.52%  '
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +33.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").Value = '2.071.53'
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0689'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.70%  '
